# Additional Address Details.xlsx — "Add files via upload" edit
# Fills in the Additional Address Grid with the applicant's address history
# and restyles the table (wrap/top-aligned text, short-date columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-style the title (row 3) and header (row 4) rows: vertical=top +
#    wrapText instead of vertical=center.
# ---------------------------------------------------------------------
$ws.Range("D3:H3").VerticalAlignment = -4160   # xlTop
$ws.Range("D3:H3").WrapText = $true

$ws.Range("D4:H4").VerticalAlignment = -4160   # xlTop
$ws.Range("D4:H4").WrapText = $true

# ---------------------------------------------------------------------
# 2) Address history data, rows 5-9 (D:H = Name, Date From, Date To,
#    Residential Address, Country)
# ---------------------------------------------------------------------
$name = "Nikita Ramesh Patil"
$country = "India"

$addr1 = "D-43/15, N-12, Swami vivekanand Nagar, HUDCO, Aurangabad`nState:Maharashtra`nPIN:431003"
$addr2 = "Survey 193/3, suvarna Building, Shankar kalate Nagar, Wakad, Pune`nState: Maharashtra`nPIN: 411057"
$addr3 = "TCS Peepul Park Rd, Technopark Campus, Thiruvananthapuram`nState: Kerala`nPIN: 695581"
$addr4 = "Adarsh Nagar Rd, Jafrabad`nState:Maharashtra`nPIN:431206"

$ws.Range("D5").Value = $name
$ws.Range("E5").Value = (Get-Date -Year 2020 -Month 5 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F5").Value = "TILL NOW"
$ws.Range("G5").Value = $addr1
$ws.Range("H5").Value = $country

$ws.Range("D6").Value = $name
$ws.Range("E6").Value = (Get-Date -Year 2016 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F6").Value = (Get-Date -Year 2020 -Month 5 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G6").Value = $addr2
$ws.Range("H6").Value = $country

$ws.Range("D7").Value = $name
$ws.Range("E7").Value = (Get-Date -Year 2015 -Month 10 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F7").Value = (Get-Date -Year 2016 -Month 2 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G7").Value = $addr3
$ws.Range("H7").Value = $country

$ws.Range("D8").Value = $name
$ws.Range("E8").Value = (Get-Date -Year 2011 -Month 7 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F8").Value = (Get-Date -Year 2015 -Month 10 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G8").Value = $addr1
$ws.Range("H8").Value = $country

$ws.Range("D9").Value = $name
$ws.Range("E9").Value = (Get-Date -Year 2004 -Month 7 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F9").Value = (Get-Date -Year 2011 -Month 7 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G9").Value = $addr4
$ws.Range("H9").Value = $country

# ---------------------------------------------------------------------
# 3) Formatting for the data rows
#    - Name / Address / Country columns: top-aligned, wrapped text
#    - Date columns: top-aligned, wrapped, short-date number format
#    - Row heights: 75pt for rows 5-8, 45pt for row 9
# ---------------------------------------------------------------------
$ws.Range("D5:D9").WrapText = $true
$ws.Range("D5:D9").VerticalAlignment = -4160

$ws.Range("G5:H9").WrapText = $true
$ws.Range("G5:H9").VerticalAlignment = -4160

$ws.Range("E5:F9").WrapText = $true
$ws.Range("E5:F9").VerticalAlignment = -4160
$ws.Range("E5:F9").NumberFormat = "mm-dd-yy"
$ws.Range("F5").NumberFormat = "General"

$ws.Rows("5:8").RowHeight = 75
$ws.Rows("9").RowHeight = 45

# ---------------------------------------------------------------------
# 4) Rows 10-13: keep blank, but now wrap-text/top-aligned like the data
#    rows above them. Row 14: blank, top-aligned only (no wrap).
# ---------------------------------------------------------------------
$ws.Range("D10:H13").WrapText = $true
$ws.Range("D10:H13").VerticalAlignment = -4160

$ws.Range("D14:H14").VerticalAlignment = -4160
